# Applies an automatic-update permutation of rows 27-32 on sheet "Artfynd".
# Observation records' identifying columns (Id, Taxonsorteringsordning, Rodlistade,
# TaxonId, Artnamn, Vetenskapligt namn, Auktor, Ost, Nord) -- and, for the
# "Vagbandad barkbock" record, its Enhet/Alder-Stadium/Kon/Aktivitet/Metod/
# Bestamningsmetod cells -- are reassigned across rows 27-32 to match the
# upstream dataset's row order. All other columns (Lokalnamn, Lan, Kommun, ...)
# are identical across these rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(27, 1).Value = 112395267
$ws.Cells.Item(27, 2).Value = 77650
$ws.Cells.Item(27, 4).Value = 'NT'
$ws.Cells.Item(27, 5).Value = 6425
$ws.Cells.Item(27, 6).Value = 'Garnlav'
$ws.Cells.Item(27, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(27, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(27, 17).Value = 331734
$ws.Cells.Item(27, 18).Value = 6626659
$ws.Cells.Item(28, 1).Value = 112395265
$ws.Cells.Item(28, 2).Value = 5135
$ws.Cells.Item(28, 4).Value = 'LC'
$ws.Cells.Item(28, 5).Value = 105930
$ws.Cells.Item(28, 6).Value = 'Vågbandad barkbock'
$ws.Cells.Item(28, 7).Value = 'Semanotus undatus'
$ws.Cells.Item(28, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(28, 10).Style = 'Normal'
$ws.Cells.Item(28, 11).Style = 'Normal'
$ws.Cells.Item(28, 12).Style = 'Normal'
$ws.Cells.Item(28, 13).Value = 'äldre gnagspår'
$ws.Cells.Item(28, 14).Style = 'Normal'
$ws.Cells.Item(28, 17).Value = 331818
$ws.Cells.Item(28, 18).Value = 6626574
$ws.Cells.Item(28, 32).Style = 'Normal'
$ws.Cells.Item(29, 1).Value = 112395263
$ws.Cells.Item(29, 2).Value = 99874
$ws.Cells.Item(29, 4).Value = 'LC'
$ws.Cells.Item(29, 5).Value = 221235
$ws.Cells.Item(29, 6).Value = 'Vårärt'
$ws.Cells.Item(29, 7).Value = 'Lathyrus vernus'
$ws.Cells.Item(29, 8).Value = '(L.) Bernh.'
$ws.Cells.Item(29, 17).Value = 331849
$ws.Cells.Item(29, 18).Value = 6626616
$ws.Cells.Item(30, 1).Value = 112395269
$ws.Cells.Item(30, 2).Value = 73772
$ws.Cells.Item(30, 4).Value = 'LC'
$ws.Cells.Item(30, 5).Value = 6426
$ws.Cells.Item(30, 6).Value = 'Kattfotslav'
$ws.Cells.Item(30, 7).Value = 'Felipes leucopellaeus'
$ws.Cells.Item(30, 8).Value = '(Ach.) Frisch & G.Thor'
$ws.Cells.Item(30, 10).Value = ''
$ws.Cells.Item(30, 11).Value = ''
$ws.Cells.Item(30, 12).Value = ''
$ws.Cells.Item(30, 13).Value = ''
$ws.Cells.Item(30, 14).Value = ''
$ws.Cells.Item(30, 17).Value = 331242
$ws.Cells.Item(30, 18).Value = 6626564
$ws.Cells.Item(30, 32).Value = ''
$ws.Cells.Item(31, 1).Value = 112395271
$ws.Cells.Item(31, 2).Value = 90814
$ws.Cells.Item(31, 4).Value = 'LC'
$ws.Cells.Item(31, 5).Value = 4364
$ws.Cells.Item(31, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(31, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(31, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(31, 17).Value = 331191
$ws.Cells.Item(31, 18).Value = 6626654
$ws.Cells.Item(32, 1).Value = 112395266
$ws.Cells.Item(32, 2).Value = 77650
$ws.Cells.Item(32, 4).Value = 'NT'
$ws.Cells.Item(32, 5).Value = 6425
$ws.Cells.Item(32, 6).Value = 'Garnlav'
$ws.Cells.Item(32, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(32, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(32, 17).Value = 331766
$ws.Cells.Item(32, 18).Value = 6626669
